$d = $word.ActiveDocument

function Find-ParaIndex($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text.TrimEnd("`r", "`a") -eq $text) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------
# 1) Remove the "Remove top of activeplan..." / "dateComponent.day..."
#    paragraphs (plus the blank line that trails them), collapsing the
#    two blank paragraphs above "@@@@" down to one.
# ---------------------------------------------------------------
$iRemoveTop = Find-ParaIndex $d "Remove top of activeplan once past time"
$iDateComp  = Find-ParaIndex $d "dateComponent.day = i-1 (depending on removal)"

$startPara = $d.Paragraphs($iRemoveTop - 1)   # the blank paragraph right before "Remove top of..."
$endPara   = $d.Paragraphs($iDateComp + 1)    # the blank paragraph right after "dateComponent...."
$r = $d.Range($startPara.Range.Start, $endPara.Range.End)
$r.Delete()

# ---------------------------------------------------------------
# 2) Remove "Lock orientation" and "Dark mode lock?" paragraphs.
# ---------------------------------------------------------------
$iLock = Find-ParaIndex $d "Lock orientation"
$iDark = Find-ParaIndex $d "Dark mode lock?"
$r2 = $d.Range($d.Paragraphs($iLock).Range.Start, $d.Paragraphs($iDark).Range.End)
$r2.Delete()

# ---------------------------------------------------------------
# 3) After the bookmark paragraph ("#A7C8CA"), insert three new notes.
# ---------------------------------------------------------------
$iBookmark = Find-ParaIndex $d "#A7C8CA"
$cur = $d.Paragraphs($iBookmark)

$cur.Range.InsertParagraphAfter()
$p = $d.Paragraphs($cur.Index + 1)
$p.Range.Text = "Logo (in assets folder App Icon at top)"
$cur = $p

$cur.Range.InsertParagraphAfter()
$p = $d.Paragraphs($cur.Index + 1)
$p.Range.Text = "Splash screen when loading (composites)"
$cur = $p

$cur.Range.InsertParagraphAfter()
$p = $d.Paragraphs($cur.Index + 1)
$p.Range.Text = "Make the list selections coloured too – if cant do that can make it small dropdown instead of full screen"
$cur = $p

# ---------------------------------------------------------------
# 4) Insert "mother does web design..." right before the "I also tested
#    with my brother..." paragraph (i.e. right after "for testing have
#    testing phases..."). Anchoring the insert on the *preceding*
#    paragraph (and using InsertParagraphAfter) avoids relying on a
#    paragraph object's .Index after it has been shifted by a mutation.
# ---------------------------------------------------------------
$iForTesting = Find-ParaIndex $d "for testing have testing phases at design, development and finalised"
$forTesting = $d.Paragraphs($iForTesting)
$forTesting.Range.InsertParagraphAfter()
$newP = $d.Paragraphs($forTesting.Index + 1)
$newP.Range.Text = "mother does web design and likes minimalistic because design often gets cluttered so my app doesn’t include spam"

# Re-resolve the "brother" paragraph fresh (it now sits right after $newP).
$brother = $d.Paragraphs($newP.Index + 1)

# ---------------------------------------------------------------
# 5) After "I also tested with my brother...", append: a blank line,
#    five new task notes for the video, then two trailing blank lines.
# ---------------------------------------------------------------
$videoTexts = @("", "For video", "Validation alert", "Notifications", "Remove days dynamically from list", "Show that changing views does not reset the form entry", "", "")

$cur = $brother
foreach ($t in $videoTexts) {
    $cur.Range.InsertParagraphAfter()
    $p = $d.Paragraphs($cur.Index + 1)
    if ($t -ne "") {
        $p.Range.Text = $t
    }
    $cur = $p
}
